# Update main GSC export data: the 2025-11-04 row has dropped out of the
# export, so remove its entire row from the "Chart" data sheet and let
# every subsequent row shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the 2025-11-04 entry (A2/B2/C2) - delete it outright so the
# rest of the table (and the shared date index) shifts up by one row.
$ws.Rows.Item(2).Delete()
